# Update "想去人数" (column F) values to match the newly generated data
# output (commit: "Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 9
$ws1.Range("F3").Value = 531
$ws1.Range("F5").Value = 503
$ws1.Range("F6").Value = 950
$ws1.Range("F7").Value = 188
$ws1.Range("F9").Value = 1016
$ws1.Range("F10").Value = 808
$ws1.Range("F11").Value = 234
$ws1.Range("F12").Value = 57
$ws1.Range("F14").Value = 813
$ws1.Range("F15").Value = 274
$ws1.Range("F16").Value = 579
$ws1.Range("F17").Value = 499
$ws1.Range("F18").Value = 1325
$ws1.Range("F20").Value = 852
$ws1.Range("F21").Value = 1169
$ws1.Range("F22").Value = 2853
$ws1.Range("F23").Value = 1390
$ws1.Range("F24").Value = 691
$ws1.Range("F26").Value = 1268
$ws1.Range("F27").Value = 59
$ws1.Range("F28").Value = 1005
$ws1.Range("F29").Value = 353
$ws1.Range("F30").Value = 3035
$ws1.Range("F31").Value = 587
$ws1.Range("F33").Value = 1386

# Sheet: 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F14").Value = 22

# Sheet: 本地生活 (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 732

# Sheet: 全部类型 (All Types - aggregated view)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 9
$ws4.Range("F3").Value = 732
$ws4.Range("F4").Value = 531
$ws4.Range("F7").Value = 503
$ws4.Range("F12").Value = 950
$ws4.Range("F13").Value = 188
$ws4.Range("F16").Value = 1016
$ws4.Range("F17").Value = 808
$ws4.Range("F18").Value = 234
$ws4.Range("F20").Value = 57
$ws4.Range("F26").Value = 813
$ws4.Range("F27").Value = 274
$ws4.Range("F28").Value = 579
$ws4.Range("F29").Value = 499
$ws4.Range("F30").Value = 1325
$ws4.Range("F32").Value = 852
$ws4.Range("F33").Value = 1169
$ws4.Range("F34").Value = 2853
$ws4.Range("F35").Value = 1390
$ws4.Range("F36").Value = 691
$ws4.Range("F38").Value = 1268
$ws4.Range("F39").Value = 59
$ws4.Range("F41").Value = 22
$ws4.Range("F42").Value = 1005
$ws4.Range("F43").Value = 353
$ws4.Range("F44").Value = 3035
$ws4.Range("F45").Value = 587
$ws4.Range("F47").Value = 1386
